# Weekly price-sheet update: a new observation (for "Ají", Vega Monumental
# Concepción) is inserted as the new row 45, pushing the previously-existing
# rows 45-71 down to 46-72 (dimension grows from A1:R71 to A1:R72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 45 - shifts old rows 45..71 down to 46..72.
$ws.Rows(45).Insert()

# Populate the newly-inserted row 45 with this week's data.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44603
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112021
$ws.Range("G45").Value = "Ají"
$ws.Range("H45").Value = "Americana (o)"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 80
$ws.Range("K45").Value = 23000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 24250
$ws.Range("N45").Value = "$/caja 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 970
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
